$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) pairs for column F ("想去人数")
$updates = @{
    "展览" = @{
        2  = 313
        4  = 3195
        8  = 7481
        9  = 73
        11 = 1213
        15 = 1049
        20 = 5866
        21 = 2326
        22 = 4054
        24 = 222
        25 = 246
        35 = 326
        37 = 942
        39 = 71
        41 = 240
    }
    "演出" = @{
        11 = 59
        18 = 121
        25 = 3437
        26 = 3437
        28 = 43
    }
    "本地生活" = @{
        4  = 575
        6  = 1906
        8  = 2967
        10 = 1197
        12 = 493
        13 = 1936
        14 = 8497
        15 = 704
    }
    "全部类型" = @{
        2  = 575
        3  = 313
        5  = 3195
        7  = 1906
        9  = 2967
        11 = 1197
        14 = 493
        16 = 1213
        23 = 1049
        24 = 59
        30 = 5866
        31 = 2326
        32 = 4054
        34 = 222
        35 = 246
        42 = 326
        46 = 71
        47 = 3437
        49 = 43
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
